$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.399.93"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "1.667.88"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'312.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.3948"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("D8").Value = "'0.3929"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'52.08"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.00%  "
$ws.Range("D10").Value = "'1.392"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.87%  "
$ws.Range("D11").Value = "'1.002"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "'0.08559"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("D13").Value = "'24.47"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.72%  "
$ws.Range("D14").Value = "'7.278"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").Value = "'7.977"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.60%  "
$ws.Range("D16").Value = "'0.00001336"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.97%  "
$ws.Range("D17").Value = "1.666.60"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "'94.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "'0.07043"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.56%  "
$ws.Range("D20").Value = "'20.58"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "'6.989"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'13.74"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("D24").Value = "24.412.96"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").Value = "'2.496"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.23%  "
$ws.Range("D26").Value = "'3.077"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +14.69%  "
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'157.15"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'142.77"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").Value = "'5.426"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").Value = "'7.960"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -9.03%  "
$ws.Range("D32").Value = "'2.551"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.42%  "
$ws.Range("D33").Value = "1.849.09"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").Value = "'1.064"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +12.56%  "
$ws.Range("D35").Value = "'0.03106"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.14%  "
$ws.Range("D36").Value = "'0.08251"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.60%  "
$ws.Range("D37").Value = "'6.901"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'11.10"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +12.96%  "
$ws.Range("D39").Value = "'0.2761"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("D40").Value = "'0.09262"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "'0.7691"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("D42").Value = "'13.70"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.77%  "
$ws.Range("D43").Value = "'1.446"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("D45").Value = "'0.7092"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.65%  "
$ws.Range("D46").Value = "'2.543"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").Value = "'4.123"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").Value = "'0.9995"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "'0.08443"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").Value = "'136.74"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("D51").Value = "'1.267"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.37%  "
